# Fruta / hortaliza, semanal
# Insert a new weekly observation row for "Poroto granado" (Feria Lagunitas
# de Puerto Montt) right after the existing row 53, pushing the remaining
# historical rows (old 54-73) down by one and growing the sheet from
# A1:R73 to A1:R74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 54..73 down to 55..74, opening up a blank row 54.
$ws.Rows.Item(54).Insert()

# Populate the newly opened row 54 with the new weekly record.
$ws.Range("A54").Value = 4
$ws.Range("B54").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C54").Value = "Los Lagos"
$ws.Range("D54").Value = 45006
$ws.Range("E54").Value = 10
$ws.Range("F54").Value = 100112030
$ws.Range("G54").Value = "Poroto granado"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 60
$ws.Range("K54").Value = 40000
$ws.Range("L54").Value = 40000
$ws.Range("M54").Value = 40000
$ws.Range("N54").Value = "`$/saco 25 kilos"
$ws.Range("O54").Value = "Región del Maule"
$ws.Range("P54").Value = 1600
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"
